# Auto-generated edit script for cryptos.xlsx update
# Commit: Updated cryptos list on Sat May  6 17:56:45 UTC 2023 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values to write: Coin (B), Link (C), Price (D), Volume 1h (E)
$updates = [ordered]@{
    'D2' = '28.868.15'
    'E2' = '  -2.64%  '
    'D3' = '1.889.32'
    'E3' = '  -5.59%  '
    'D4' = '1.003'
    'E4' = '  -0.12%  '
    'D5' = '323.17'
    'E5' = '  -2.05%  '
    'E6' = '  -0.33%  '
    'D7' = '0.4585'
    'E7' = '  -1.98%  '
    'B9' = 'OKB'
    'C9' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'D9' = '45.61'
    'E9' = '  -2.37%  '
    'B10' = 'Dogecoin'
    'C10' = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
    'D10' = '0.07717'
    'E10' = '  -3.26%  '
    'B11' = 'Polygon'
    'C11' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'D11' = '0.9654'
    'E11' = '  -4.10%  '
    'B12' = 'Solana'
    'C12' = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
    'D12' = '22.03'
    'E12' = '  -3.20%  '
    'B13' = 'WrappedEther'
    'C13' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D13' = '1.886.62'
    'E13' = '  -7.04%  '
    'B14' = 'Chainlink'
    'C14' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'D14' = '6.945'
    'E14' = '  -4.27%  '
    'B15' = 'Polkadot'
    'C15' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'D15' = '5.657'
    'E15' = '  -3.85%  '
    'B16' = 'TRON'
    'C16' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'D16' = '0.07039'
    'E16' = '  -1.93%  '
    'B17' = 'BinanceUSD'
    'C17' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'D17' = '1.004'
    'E17' = '  -0.17%  '
    'B18' = 'Litecoin'
    'C18' = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    'D18' = '83.25'
    'E18' = '  -6.54%  '
    'B19' = 'ShibaInu'
    'C19' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'D19' = '0.000009503'
    'E19' = '  -4.99%  '
    'B20' = 'Avalanche'
    'C20' = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    'D20' = '16.65'
    'E20' = '  -3.36%  '
    'B21' = 'Dai'
    'C21' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'D21' = '1.002'
    'E21' = '  -0.27%  '
    'B22' = 'WrappedBTC'
    'C22' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    'D22' = '28.824.30'
    'E22' = '  -2.98%  '
    'B23' = 'Uniswap'
    'C23' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D23' = '5.307'
    'E23' = '  -4.32%  '
    'B24' = 'Cosmos'
    'C24' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D24' = '10.87'
    'E24' = '  -3.71%  '
    'B25' = 'WrappedliquidstakedEther2.0'
    'C25' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'D25' = '2.113.40'
    'E25' = '  -6.77%  '
    'B26' = 'Toncoin'
    'C26' = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    'D26' = '2.078'
    'E26' = '  -3.07%  '
    'B27' = 'Monero'
    'C27' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D27' = '156.18'
    'E27' = '  -1.51%  '
    'B28' = 'EthereumClassic'
    'C28' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D28' = '19.00'
    'E28' = '  -3.63%  '
    'B29' = 'InternetComputer(DFINITY)'
    'C29' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D29' = '5.585'
    'E29' = '  -6.99%  '
    'B30' = 'BitcoinCash'
    'C30' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'D30' = '117.15'
    'E30' = '  -3.02%  '
    'B31' = 'LidoDAOToken'
    'C31' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'D31' = '1.808'
    'E31' = '  -8.07%  '
    'B32' = 'Stellar'
    'C32' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D32' = '0.09270'
    'E32' = '  -1.99%  '
    'B33' = 'ImmutableX'
    'C33' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D33' = '0.8492'
    'E33' = '  -5.26%  '
    'B34' = 'Filecoin'
    'C34' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D34' = '5.073'
    'E34' = '  -4.53%  '
    'B35' = 'ARBITRUM'
    'C35' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D35' = '1.224'
    'E35' = '  -9.41%  '
    'B36' = 'HuobiToken'
    'C36' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D36' = '3.014'
    'E36' = '  -5.45%  '
    'B37' = 'Hedera'
    'C37' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D37' = '0.05672'
    'E37' = '  -3.13%  '
    'B38' = 'TrustWalletToken'
    'C38' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D38' = '1.138'
    'E38' = '  -3.84%  '
    'B39' = 'Frax'
    'C39' = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
    'D39' = '1.002'
    'E39' = '  -0.10%  '
    'B40' = 'VeChain'
    'C40' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D40' = '0.02035'
    'E40' = '  -4.48%  '
    'B41' = 'FraxShare'
    'C41' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D41' = '7.398'
    'E41' = '  -6.89%  '
    'B42' = 'TheSandbox'
    'C42' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D42' = '0.5481'
    'E42' = '  -5.23%  '
    'B43' = 'Algorand'
    'C43' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D43' = '0.1747'
    'E43' = '  -4.25%  '
    'B44' = 'PEPE'
    'C44' = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
    'D44' = '0.000002870'
    'E44' = '  -23.98%  '
    'B45' = 'Aptos'
    'C45' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D45' = '9.161'
    'E45' = '  -7.62%  '
    'B46' = 'MXToken'
    'C46' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D46' = '2.703'
    'E46' = '  +1.75%  '
    'B47' = 'Decentraland'
    'C47' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    'D47' = '0.5151'
    'E47' = '  -4.81%  '
    'B48' = 'EnergySwap'
    'C48' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D48' = '11.19'
    'E48' = '  -8.30%  '
    'B49' = 'Cronos'
    'C49' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D49' = '0.06806'
    'E49' = '  -2.99%  '
    'B50' = 'RenderToken'
    'C50' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D50' = '2.065'
    'E50' = '  -5.18%  '
    'B51' = 'Quant'
    'C51' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D51' = '111.38'
    'E51' = '  -3.00%  '
}

# Columns whose new values are plain decimal-looking strings (e.g. "1.003")
# must have their number format forced to Text first, otherwise Excel
# auto-converts them into numeric values (losing formatting / precision).
$textForceColumns = @("D")

foreach ($ref in $updates.Keys) {
    $col = $ref -replace "[0-9]+$", ""
    $cell = $ws.Range($ref)
    if ($textForceColumns -contains $col) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $updates[$ref]
}

Write-Host "Applied $($updates.Count) cell updates"
